$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.533.88'
$ws.Range("E2").Value = '  +2.00%  '
$ws.Range("D3").Value = '3.557.30'
$ws.Range("E3").Value = '  +2.12%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '620.59'
$ws.Range("E5").Value = '  +2.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.22'
$ws.Range("E6").Value = '  +5.58%  '
$ws.Range("D7").Value = '3.550.24'
$ws.Range("E7").Value = '  +1.94%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.490'
$ws.Range("E9").Value = '  +2.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.147'
$ws.Range("E10").Value = '  +6.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.31'
$ws.Range("E11").Value = '  +5.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.439'
$ws.Range("E12").Value = '  +4.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000224'
$ws.Range("E13").Value = '  +4.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.30'
$ws.Range("E14").Value = '  +6.27%  '
$ws.Range("D15").Value = '4.166.53'
$ws.Range("E15").Value = '  +2.30%  '
$ws.Range("D16").Value = '69.093.18'
$ws.Range("E16").Value = '  +2.96%  '
$ws.Range("D17").Value = '3.570.73'
$ws.Range("E17").Value = '  +2.34%  '
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.81'
$ws.Range("E19").Value = '  +6.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.97'
$ws.Range("E20").Value = '  +6.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.08'
$ws.Range("E21").Value = '  +11.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '458.55'
$ws.Range("E22").Value = '  +2.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.642'
$ws.Range("E23").Value = '  +3.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.59'
$ws.Range("E24").Value = '  +2.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000131'
$ws.Range("E25").Value = '  +4.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.54'
$ws.Range("E26").Value = '  +4.48%  '
$ws.Range("D27").Value = '3.713.47'
$ws.Range("E27").Value = '  +2.51%  '
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.09'
$ws.Range("E29").Value = '  +9.98%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.57'
$ws.Range("E30").Value = '  +2.04%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.70'
$ws.Range("E31").Value = '  +8.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.170'
$ws.Range("E32").Value = '  +4.06%  '
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.23%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.47'
$ws.Range("E34").Value = '  +6.11%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '26.22'
$ws.Range("E35").Value = '  +2.09%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.93'
$ws.Range("E36").Value = '  +4.64%  '
$ws.Range("D37").Value = '3.554.09'
$ws.Range("E37").Value = '  +2.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.30'
$ws.Range("E38").Value = '  +4.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.39'
$ws.Range("E39").Value = '  +8.80%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '178.52'
$ws.Range("E41").Value = '  +4.79%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.41%  '
$ws.Range("B43").Value = 'Hedera'
$ws.Range("C43").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0917'
$ws.Range("E43").Value = '  +5.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.60'
$ws.Range("E44").Value = '  +3.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '30.65'
$ws.Range("E45").Value = '  +15.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.898'
$ws.Range("E46").Value = '  +1.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.34'
$ws.Range("E47").Value = '  +7.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '45.72'
$ws.Range("E48").Value = '  +0.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.68'
$ws.Range("E49").Value = '  +6.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.80'
$ws.Range("E50").Value = '  +3.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.264'
$ws.Range("E51").Value = '  +8.74%  '
